# "#5: fund, bonds, otherbonds, antique done"
#
#  - Finish the "基金受益憑證" (fund) sheet: row 1 was accidentally filled
#    with a copy of row 2's data instead of real column headers; replace it
#    with the proper header labels and extend row 2 with the common
#    trailing metadata columns (property_category .. index) used by every
#    other finished sheet in this workbook.
#  - The "其他有價證券" (other securities) sheet isn't done yet this round,
#    so it is removed again; "保險" (insurance) simply shifts up to take
#    its place.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Fix up the "基金受益憑證" (fund) sheet.
# ---------------------------------------------------------------------
$fund  = $wb.Worksheets.Item("基金受益憑證")
$stock = $wb.Worksheets.Item("股票")   # already-finished sheet to borrow a
                                        # literal "2012-04-23" text cell from

# -- Row 1: turn the bogus duplicate data into real headers --------------
# Copy the existing bold/bordered header style from B1 onto the new header
# cells first so they match the look of the rest of the row.
$fund.Range("B1:H1").Copy() | Out-Null
$fund.Range("I1:O1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$fund.Cells.Item(1,2).Value  = "name"
$fund.Cells.Item(1,3).Value  = "owner"
$fund.Cells.Item(1,4).Value  = "dealer"
$fund.Cells.Item(1,5).Value  = "quantity"
$fund.Cells.Item(1,6).Value  = "face_value"
$fund.Cells.Item(1,7).Value  = "currency"
$fund.Cells.Item(1,8).Value  = "total"
$fund.Cells.Item(1,9).Value  = "property_category"
$fund.Cells.Item(1,10).Value = "category"
$fund.Cells.Item(1,11).Value = "date"
$fund.Cells.Item(1,12).Value = "legislator_name"
$fund.Cells.Item(1,13).Value = "legislator_id"
$fund.Cells.Item(1,14).Value = "source_file"
$fund.Cells.Item(1,15).Value = "index"

# -- Row 2: keep name/owner/dealer/quantity/face_value/currency/total and
#    append the shared trailing metadata columns. ------------------------
$fund.Cells.Item(2,9).Value  = "fund"
$fund.Cells.Item(2,10).Value = "normal"

# "2012-04-23" looks like a date, so typing it in would get reinterpreted
# as a date serial. Instead copy the value straight from another finished
# sheet's cell that already holds the same literal text, which keeps it a
# plain text cell.
$stock.Range("J2").Copy() | Out-Null
$fund.Range("K2").PasteSpecial(-4163) | Out-Null   # xlPasteValues

$fund.Cells.Item(2,12).Value = "許忠信"
$fund.Cells.Item(2,13).Value = 1749
$fund.Cells.Item(2,14).Value = "tmpa22c1"
$fund.Cells.Item(2,15).Value = 96

# ---------------------------------------------------------------------
# 2) Drop the still-unfinished "其他有價證券" sheet.
# ---------------------------------------------------------------------
$other = $wb.Worksheets.Item("其他有價證券")
$other.Delete()
